# Daily attendance processing - reorder "Recorded By" (column G) entries so
# that any "System"/"system" token(s) are moved to the front of the
# comma-separated list, while the remaining tokens keep their original
# relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = 7
    $val = $cell.Value()

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "

    $sysParts = @()
    $otherParts = @()
    foreach ($p in $parts) {
        if ($p.ToLower() -eq "system") {
            $sysParts += $p
        } else {
            $otherParts += $p
        }
    }

    if ($sysParts.Count -gt 0 -and $otherParts.Count -gt 0) {
        $newVal = ($sysParts + $otherParts) -join ", "
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}
